# Automated update of BRVM recommendations data (Recommandations + Top_YTD sheets)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations" updates ---
$ws1.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 8
$ws1.Cells.Item(2, 4).Value = 3366.19
$ws1.Cells.Item(2, 5).Value = 108.93
$ws1.Cells.Item(3, 3).Value = 3
$ws1.Cells.Item(3, 4).Value = 2835
$ws1.Cells.Item(3, 5).Value = 935
$ws1.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 3
$ws1.Cells.Item(4, 4).Value = 2830
$ws1.Cells.Item(4, 5).Value = 940
$ws1.Cells.Item(6, 4).Value = 2670
$ws1.Cells.Item(6, 5).Value = 690
$ws1.Cells.Item(7, 4).Value = 2557.71
$ws1.Cells.Item(7, 5).Value = 649.41
$ws1.Cells.Item(8, 4).Value = 1620.87
$ws1.Cells.Item(8, 5).Value = 413.43
$ws1.Cells.Item(9, 4).Value = 1450.28
$ws1.Cells.Item(9, 5).Value = 367.16
$ws1.Cells.Item(10, 4).Value = 1359.27
$ws1.Cells.Item(10, 5).Value = 344.57
$ws1.Cells.Item(12, 4).Value = 569.5599999999999
$ws1.Cells.Item(12, 5).Value = 143.57
$ws1.Cells.Item(13, 4).Value = 546.55
$ws1.Cells.Item(13, 5).Value = 136.52
$ws1.Cells.Item(14, 4).Value = 535
$ws1.Cells.Item(14, 5).Value = 133.97
$ws1.Cells.Item(15, 4).Value = 525.78
$ws1.Cells.Item(15, 5).Value = 131.66
$ws1.Cells.Item(16, 4).Value = 515.89
$ws1.Cells.Item(16, 5).Value = 134
$ws1.Cells.Item(17, 4).Value = 429.24
$ws1.Cells.Item(17, 5).Value = 107.48
$ws1.Cells.Item(18, 4).Value = 381.71
$ws1.Cells.Item(18, 5).Value = 94.98999999999999
$ws1.Cells.Item(19, 1).Value = "BRVM - CONSOMMATION DE BASE            (**)"
$ws1.Cells.Item(19, 2).Value = 0
$ws1.Cells.Item(19, 3).Value = 2
$ws1.Cells.Item(19, 4).Value = 371.92
$ws1.Cells.Item(19, 5).Value = 185.85
$ws1.Cells.Item(22, 1).Value = "BRVM - INDUSTRIE                       (**)"
$ws1.Cells.Item(22, 4).Value = 214.57
$ws1.Cells.Item(22, 5).Value = 214.57
$ws1.Cells.Item(23, 1).Value = "BRVM-PRINCIPAL                          (**)"
$ws1.Cells.Item(23, 4).Value = 192.26
$ws1.Cells.Item(23, 5).Value = 192.26
$ws1.Cells.Item(24, 1).Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Cells.Item(24, 4).Value = 192.05
$ws1.Cells.Item(24, 5).Value = 192.05
$ws1.Cells.Item(25, 1).Value = "BRVM-PRINCIPAL                       (**)"
$ws1.Cells.Item(25, 4).Value = 191.44
$ws1.Cells.Item(25, 5).Value = 191.44
$ws1.Cells.Item(27, 1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(27, 2).Value = 3
$ws1.Cells.Item(27, 3).Value = 0
$ws1.Cells.Item(27, 4).Value = 21.64
$ws1.Cells.Item(27, 5).Value = 7.49
$ws1.Cells.Item(27, 6).Value = "🟢 Achat"
$ws1.Cells.Item(27, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(30, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(30, 2).Value = 1
$ws1.Cells.Item(30, 3).Value = 0
$ws1.Cells.Item(30, 4).Value = 7.49
$ws1.Cells.Item(30, 5).Value = 7.49
$ws1.Cells.Item(31, 1).Value = "BANK OF AFRICA CI (BOAC)"
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 0
$ws1.Cells.Item(31, 4).Value = 7.44
$ws1.Cells.Item(31, 5).Value = 7.44
$ws1.Cells.Item(32, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(32, 4).Value = 5.76
$ws1.Cells.Item(32, 5).Value = 5.76
$ws1.Cells.Item(33, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(33, 4).Value = 4.87
$ws1.Cells.Item(33, 5).Value = 4.87
$ws1.Cells.Item(34, 1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(34, 2).Value = 1
$ws1.Cells.Item(34, 3).Value = 1
$ws1.Cells.Item(34, 4).Value = 4.5
$ws1.Cells.Item(34, 5).Value = -2.96
$ws1.Cells.Item(34, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(35, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(35, 2).Value = 1
$ws1.Cells.Item(35, 3).Value = 0
$ws1.Cells.Item(35, 4).Value = 3.9
$ws1.Cells.Item(35, 5).Value = 3.9
$ws1.Cells.Item(35, 7).Value = "➖ Neutre"
$ws1.Cells.Item(36, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(36, 4).Value = 3.7
$ws1.Cells.Item(36, 5).Value = 3.7
$ws1.Cells.Item(37, 1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(37, 4).Value = 3.19
$ws1.Cells.Item(37, 5).Value = 3.19
$ws1.Cells.Item(38, 1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(38, 4).Value = 2.61
$ws1.Cells.Item(38, 5).Value = 5.26
$ws1.Cells.Item(39, 1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(39, 4).Value = 1.31
$ws1.Cells.Item(39, 5).Value = 7.48
$ws1.Cells.Item(40, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(40, 2).Value = 1
$ws1.Cells.Item(40, 3).Value = 1
$ws1.Cells.Item(40, 4).Value = 0.32
$ws1.Cells.Item(40, 5).Value = 5.88
$ws1.Cells.Item(40, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(41, 1).Value = "TOTAL"
$ws1.Cells.Item(41, 2).Value = 0
$ws1.Cells.Item(41, 3).Value = 4
$ws1.Cells.Item(41, 4).Value = 0
$ws1.Cells.Item(41, 5).Value = 0
$ws1.Cells.Item(41, 7).Value = "➖ Neutre"
$ws1.Cells.Item(45, 1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(45, 2).Value = 0
$ws1.Cells.Item(45, 3).Value = 1
$ws1.Cells.Item(45, 4).Value = -2.03
$ws1.Cells.Item(45, 5).Value = -2.03
$ws1.Cells.Item(45, 7).Value = "➖ Neutre"
$ws1.Cells.Item(48, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(48, 4).Value = -3.7
$ws1.Cells.Item(48, 5).Value = -3.7
$ws1.Cells.Item(49, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(49, 4).Value = -4.88
$ws1.Cells.Item(49, 5).Value = -4.88
$ws1.Cells.Item(50, 1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(50, 2).Value = 0
$ws1.Cells.Item(50, 3).Value = 3
$ws1.Cells.Item(50, 4).Value = -7.7
$ws1.Cells.Item(50, 5).Value = -2.6
$ws1.Cells.Item(50, 6).Value = "🔴 Vente"
$ws1.Cells.Item(50, 7).Value = "⚠️ Risque de décrochage"
$ws1.Cells.Item(51, 4).Value = -29.84
$ws1.Cells.Item(51, 5).Value = -7.49

# --- Sheet "Top_YTD" updates ---
$ws2.Cells.Item(2, 2).Value = 9179936.130000001
$ws2.Cells.Item(3, 1).Value = "AIR LIQUIDE CI"
$ws2.Cells.Item(3, 2).Value = 394346.8
$ws2.Cells.Item(4, 1).Value = "NEI-CEDA CI"
$ws2.Cells.Item(4, 2).Value = 346676.03
$ws2.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$ws2.Cells.Item(5, 2).Value = 298716.78
$ws2.Cells.Item(6, 1).Value = "SETAO CI"
$ws2.Cells.Item(6, 2).Value = 113944.06
$ws2.Cells.Item(7, 1).Value = "CFAO MOTORS CI"
$ws2.Cells.Item(7, 2).Value = 113468
$ws2.Cells.Item(8, 2).Value = 65036.41
$ws2.Cells.Item(9, 2).Value = 45676.4
$ws2.Cells.Item(10, 2).Value = 37315.25
$ws2.Cells.Item(11, 2).Value = 3351.82
